$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")
Write-Host $ws.Name
